$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had an extra leading column A (GENE-repeat values:
# 2, 6, 7, 11) that duplicated the data already present in column F.
# Remove that stray column so columns B:F shift left to become A:E,
# matching the "all result ready to start write" cleanup.
$ws.Range("A1").EntireColumn.Delete()
